$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows for "Huynh Tu Canh" under tieu-muc 1003/1701 (original rows 5-6)
# were duplicate/erroneous "truy thu" (back-tax) lines at a 100% tax rate;
# remove them so every row below shifts up by two.
$ws.Rows("5:6").Delete()

# Renumber the STT (serial number) column for the remaining 9 data rows.
for ($i = 0; $i -le 8; $i++) {
    $ws.Cells.Item(5 + $i, 1).Value = $i + 1
}

# Column A on the rows that used to belong to "normal" (non-highlighted)
# entries lost its bordered/highlighted look after the shift; bring it back
# in line with the rest of the table by copying the format from A5.
$ws.Range("A5").Copy()
$ws.Range("A6:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect where the user's cursor ended up after the edit.
$ws.Range("C11").Select() | Out-Null
